# Fruta / hortaliza, semanal
#
# A new weekly price-report row is inserted into the "Vega Modelo de Temuco -
# Maracuyá" sheet at row 20 (pushing the existing rows 20-31 down to 21-32).
# The new row carries a fresh observation; every column except the shifted
# rows' own values is otherwise untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 20 - everything from row 20 downward
# (through the old row 31) shifts down by one, becoming rows 21-32.
$ws.Rows(20).Insert()

# Populate the newly-inserted row 20 with the new weekly observation.
$ws.Cells.Item(20, 1).Value = 10
$ws.Cells.Item(20, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(20, 3).Value = "La Araucanía"
$ws.Cells.Item(20, 4).Value = 44669
$ws.Cells.Item(20, 5).Value = 9
$ws.Cells.Item(20, 6).Value = "Fruta"
$ws.Cells.Item(20, 7).Value = 100108
$ws.Cells.Item(20, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(20, 9).Value = 100108003
$ws.Cells.Item(20, 10).Value = "Maracuyá"
$ws.Cells.Item(20, 11).Value = "Sin especificar"
$ws.Cells.Item(20, 12).Value = "Primera"
$ws.Cells.Item(20, 13).Value = 40
$ws.Cells.Item(20, 14).Value = 32000
$ws.Cells.Item(20, 15).Value = 32000
$ws.Cells.Item(20, 16).Value = 32000
$ws.Cells.Item(20, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(20, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(20, 19).Value = 1778
$ws.Cells.Item(20, 20).Value = 18
